# Apply the Silabus Non-Akademik template update:
#  - Reword the instructions cell (A1)
#  - Insert three new columns (Tahapan, Penilaian, Referensi) between
#    "Materi" and "PIC", shifting Waktu/Tempat/Sasaran right and
#    dropping the old "Catatan (opsional)" column
#  - Update the sample data row to match the new columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: instructions text ---
$ws.Range("A1").Value = "Petunjuk: Isi baris-baris silabus non-akademik. Kolom Waktu isi teks (misal: Harian/Pekanan/Bulanan) atau kosong."

# --- Row 3: header row (new layout) ---
$ws.Range("A3").Value = "Materi"
$ws.Range("B3").Value = "Tahapan (target per jenjang)"
$ws.Range("C3").Value = "Penilaian"
$ws.Range("D3").Value = "Referensi"
$ws.Range("E3").Value = "PIC"
$ws.Range("F3").Value = "Waktu"
$ws.Range("G3").Value = "Tempat"
$ws.Range("H3").Value = "Sasaran"

# Bold styling carries through A3:F3 (same as before), while the
# shifted-in G3/H3 (Tempat/Sasaran) are plain, unstyled cells.
$ws.Range("A3:F3").Font.Bold = $true
$ws.Range("G3:H3").Font.Bold = $false

# --- Row 4: sample data (new layout) ---
$ws.Range("A4").Value = "Contoh: Hiwar - Perkenalan"
$ws.Range("B4").Value = "MTs/MA (sesuaikan jenjang)"
$ws.Range("C4").Value = "Observasi + praktik"
$ws.Range("D4").Value = "Modul internal / kitab / pedoman"
$ws.Range("E4").Value = "Ust A"
$ws.Range("F4").Value = "Pekanan"
$ws.Range("G4").Value = "Aula Bahasa"
$ws.Range("H4").Value = "Kelas X"
